$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column price values are forced to text (matching source inline-string cells),
# then formats are cleared so no stray style index is introduced.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.236.06"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.89%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.487.37"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "603.11"
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "143.86"
$c.ClearFormats()
$ws.Range("E6").Value = "  -2.82%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.486.33"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("E11").Value = "  -5.11%  "

$ws.Range("E12").Value = "  -2.89%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.077.61"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.24%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "30.35"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("E15").Value = "  -5.39%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.480.60"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.36%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "66.243.87"
$c.ClearFormats()

$ws.Range("E18").Value = "  -0.32%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.81"
$c.ClearFormats()
$ws.Range("E19").Value = "  +4.29%  "

$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("E21").Value = "  -3.47%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "426.06"
$c.ClearFormats()
$ws.Range("E22").Value = "  -1.85%  "

$ws.Range("E23").Value = "  -2.58%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "77.87"
$c.ClearFormats()
$ws.Range("E24").Value = "  -2.16%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  -2.94%  "

$ws.Range("E27").Value = "  -5.93%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.ClearFormats()
$ws.Range("E28").Value = "  -3.73%  "

$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("E30").Value = "  -0.14%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.ClearFormats()
$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("E32").Value = "  -8.78%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "25.07"
$c.ClearFormats()
$ws.Range("E33").Value = "  -1.34%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.475.22"
$c.ClearFormats()
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -3.81%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.60"
$c.ClearFormats()
$ws.Range("E37").Value = "  -6.01%  "

$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("E39").Value = "  +0.01%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "170.03"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.35%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0858"
$c.ClearFormats()

$ws.Range("E42").Value = "  -4.88%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.879"
$c.ClearFormats()
$ws.Range("E43").Value = "  -1.85%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.90"
$c.ClearFormats()
$ws.Range("E44").Value = "  -8.64%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "45.42"
$c.ClearFormats()
$ws.Range("E45").Value = "  -0.90%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "26.02"
$c.ClearFormats()
$ws.Range("E46").Value = "  -10.09%  "

$ws.Range("E47").Value = "  -2.21%  "

$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("E49").Value = "  -4.59%  "

$ws.Range("E50").Value = "  -3.21%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.236"
$c.ClearFormats()
$ws.Range("E51").Value = "  -4.37%  "
